$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H6").Value = 403.16666
$ws.Range("I6").Value = 383.8
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 1151.4
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = -1039.4
$ws.Range("N6").Value = -1724
$ws.Range("H17").Value = 608.88605
$ws.Range("J17").Value = 635.49335
$ws.Range("L17").Value = 1906.48005
$ws.Range("N17").Value = -2242.48005
$ws.Range("H40").Value = 333333340
$ws.Range("I40").Value = 333333340
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 333333340
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -333333165
$ws.Range("N40").ClearContents()
$ws.Range("H43").Value = 3503.4614
$ws.Range("I43").Value = 650.5
$ws.Range("J43").Value = 4022.182
$ws.Range("K43").Value = 650.5
$ws.Range("L43").Value = 4022.182
$ws.Range("M43").Value = -581.5
$ws.Range("N43").Value = -4160.182
$ws.Range("H69").Value = 3714.2856
$ws.Range("I69").Value = 2000
$ws.Range("J69").Value = 4000
$ws.Range("K69").Value = 6000
$ws.Range("L69").Value = 12000
$ws.Range("M69").Value = -5126
$ws.Range("N69").Value = -13748
$ws.Range("H72").Value = 3714.2856
$ws.Range("I72").Value = 2000
$ws.Range("J72").Value = 4000
$ws.Range("K72").Value = 18000
$ws.Range("L72").Value = 36000
$ws.Range("M72").Value = -13632
$ws.Range("N72").Value = -44736
$ws.Range("H127").Value = 893.9474
$ws.Range("I127").Value = 476.55554
$ws.Range("J127").Value = 1269.6
$ws.Range("K127").Value = 1429.66662
$ws.Range("L127").Value = 3808.8
$ws.Range("M127").Value = 3530.33338
$ws.Range("N127").Value = -13728.8
$ws.Range("H137").Value = 25642888
$ws.Range("I137").Value = 1213.2
$ws.Range("J137").Value = 111115140
$ws.Range("K137").Value = 3639.6
$ws.Range("L137").Value = 333345420
$ws.Range("M137").Value = -1089.6
$ws.Range("N137").Value = -333350520

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H74").Value = 9200.764999999999
$ws.Range("I74").Value = 1259.1428
$ws.Range("J74").Value = 14759.9
$ws.Range("K74").Value = 1259.1428
$ws.Range("L74").Value = 14759.9
$ws.Range("M74").Value = -385.1428000000001
$ws.Range("N74").Value = -16507.9
$ws.Range("H77").Value = 9200.764999999999
$ws.Range("I77").Value = 1259.1428
$ws.Range("J77").Value = 14759.9
$ws.Range("K77").Value = 6295.714
$ws.Range("L77").Value = 73799.5
$ws.Range("M77").Value = -1927.714
$ws.Range("N77").Value = -82535.5
$ws.Range("H102").Value = 2163.75
$ws.Range("I102").Value = 2163.75
$ws.Range("K102").Value = 2163.75
$ws.Range("M102").Value = -541.75
$ws.Range("H124").Value = 9904.5
$ws.Range("J124").Value = 9904.5
$ws.Range("L124").Value = 9904.5
$ws.Range("N124").Value = -19724.5
$ws.Range("H125").Value = 38340
$ws.Range("J125").Value = 38340
$ws.Range("L125").Value = 38340
$ws.Range("N125").Value = -48180

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H10").Value = 928.5714
$ws.Range("I10").Value = 928.5714
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 928.5714
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -788.5714
$ws.Range("N10").ClearContents()
$ws.Range("H134").Value = 54174.156
$ws.Range("I134").Value = 68159.734
$ws.Range("J134").Value = 1728.25
$ws.Range("K134").Value = 204479.202
$ws.Range("L134").Value = 5184.75
$ws.Range("M134").Value = -201944.202
$ws.Range("N134").Value = -10254.75

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H10").Value = 1750.875
$ws.Range("I10").Value = 1750.875
$ws.Range("K10").Value = 1750.875
$ws.Range("M10").Value = -1611.875
$ws.Range("H31").Value = 2017.3334
$ws.Range("I31").Value = 1680
$ws.Range("J31").Value = 2312.5
$ws.Range("K31").Value = 1680
$ws.Range("L31").Value = 2312.5
$ws.Range("M31").Value = -1385
$ws.Range("N31").Value = -2902.5
$ws.Range("H34").Value = 2017.3334
$ws.Range("I34").Value = 1680
$ws.Range("J34").Value = 2312.5
$ws.Range("K34").Value = 1680
$ws.Range("L34").Value = 2312.5
$ws.Range("M34").Value = -1478
$ws.Range("N34").Value = -2716.5
$ws.Range("H58").Value = 1181.3438
$ws.Range("I58").Value = 929.6087
$ws.Range("J58").Value = 1824.6666
$ws.Range("K58").Value = 929.6087
$ws.Range("L58").Value = 1824.6666
$ws.Range("M58").Value = -726.6087
$ws.Range("N58").Value = -2230.6666
$ws.Range("H136").Value = 1181.3438
$ws.Range("I136").Value = 929.6087
$ws.Range("J136").Value = 1824.6666
$ws.Range("K136").Value = 2788.8261
$ws.Range("L136").Value = 5473.9998
$ws.Range("M136").Value = -238.8261000000002
$ws.Range("N136").Value = -10573.9998

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("I5").Value = 900
$ws.Range("K5").Value = 2700
$ws.Range("M5").Value = -2588
$ws.Range("H20").Value = 2374.75
$ws.Range("J20").Value = 2374.75
$ws.Range("L20").Value = 7124.25
$ws.Range("N20").Value = -7578.25
$ws.Range("H22").Value = 712.7143
$ws.Range("I22").Value = 496.33334
$ws.Range("J22").Value = 875
$ws.Range("K22").Value = 1489.00002
$ws.Range("L22").Value = 2625
$ws.Range("M22").Value = -1320.00002
$ws.Range("N22").Value = -2963
$ws.Range("H27").Value = 712.7143
$ws.Range("I27").Value = 496.33334
$ws.Range("J27").Value = 875
$ws.Range("K27").Value = 1489.00002
$ws.Range("L27").Value = 2625
$ws.Range("M27").Value = -1387.00002
$ws.Range("N27").Value = -2829
$ws.Range("H96").Value = 47140104
$ws.Range("J96").Value = 47140104
$ws.Range("L96").Value = 141420312
$ws.Range("N96").Value = -141424430
$ws.Range("H113").Value = 699.46155
$ws.Range("I113").Value = 795
$ws.Range("J113").Value = 682.0909
$ws.Range("K113").Value = 2385
$ws.Range("L113").Value = 2046.2727
$ws.Range("M113").Value = -215
$ws.Range("N113").Value = -6386.2727
$ws.Range("H122").Value = 20835250
$ws.Range("J122").Value = 2861
$ws.Range("L122").Value = 25749
$ws.Range("N122").Value = -30649
$ws.Range("H132").Value = 43479468
$ws.Range("J132").Value = 1590
$ws.Range("L132").Value = 14310
$ws.Range("N132").Value = -19370
$ws.Range("I135").Value = 900
$ws.Range("K135").Value = 8100
$ws.Range("M135").Value = -5565
$ws.Range("H137").Value = 42966.06
$ws.Range("I137").Value = 2674.1667
$ws.Range("J137").Value = 55053.625
$ws.Range("K137").Value = 8022.500100000001
$ws.Range("L137").Value = 165160.875
$ws.Range("M137").Value = -2922.500100000001
$ws.Range("N137").Value = -175360.875

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H126").Value = 1689.1
$ws.Range("I126").Value = 1597.8
$ws.Range("J126").Value = 1780.4
$ws.Range("K126").Value = 4793.4
$ws.Range("L126").Value = 5341.200000000001
$ws.Range("M126").Value = -2323.4
$ws.Range("N126").Value = -10281.2

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H22").Value = 402.74075
$ws.Range("I22").Value = 519
$ws.Range("J22").Value = 376.31818
$ws.Range("K22").Value = 519
$ws.Range("L22").Value = 376.31818
$ws.Range("M22").Value = -224
$ws.Range("N22").Value = -966.31818
$ws.Range("H27").Value = 402.74075
$ws.Range("I27").Value = 519
$ws.Range("J27").Value = 376.31818
$ws.Range("K27").Value = 519
$ws.Range("L27").Value = 376.31818
$ws.Range("M27").Value = -412
$ws.Range("N27").Value = -590.31818
$ws.Range("H40").Value = 1216.25
$ws.Range("I40").Value = 1269
$ws.Range("K40").Value = 1269
$ws.Range("M40").Value = -1133
$ws.Range("H100").Value = 3146.238
$ws.Range("I100").Value = 1872.091
$ws.Range("J100").Value = 4547.8
$ws.Range("K100").Value = 1872.091
$ws.Range("L100").Value = 4547.8
$ws.Range("M100").Value = -1331.091
$ws.Range("N100").Value = -5629.8

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H81").Value = 4931.25
$ws.Range("I81").Value = 900
$ws.Range("J81").Value = 5200
$ws.Range("K81").Value = 1800
$ws.Range("L81").Value = 10400
$ws.Range("M81").Value = -739
$ws.Range("N81").Value = -12522
$ws.Range("H84").Value = 4931.25
$ws.Range("I84").Value = 900
$ws.Range("J84").Value = 5200
$ws.Range("K84").Value = 9000
$ws.Range("L84").Value = 52000
$ws.Range("M84").Value = -3696
$ws.Range("N84").Value = -62608
$ws.Range("H113").Value = 579.2222
$ws.Range("I113").Value = 1050
$ws.Range("J113").Value = 444.7143
$ws.Range("K113").Value = 3150
$ws.Range("L113").Value = 1334.1429
$ws.Range("M113").Value = -980
$ws.Range("N113").Value = -5674.1429
$ws.Range("H132").Value = 2436.5454
$ws.Range("I132").Value = 2343.2432
$ws.Range("J132").Value = 2929.7144
$ws.Range("K132").Value = 7029.7296
$ws.Range("L132").Value = 8789.143199999999
$ws.Range("M132").Value = -4499.7296
$ws.Range("N132").Value = -13849.1432
